# --- Prix Spot: add column E (18-jun) ---
$wb = $excel.ActiveWorkbook
$ws0 = $wb.Worksheets.Item(1)

# Copy the header formatting from D1 into the new E1 header cell, then set value/text
$ws0.Range("D1").Copy()
$ws0.Range("E1").PasteSpecial(-4122)
$ws0.Range("E1").Value = "18-jun"

$spotValues = @(51.35, 28.31, 26.87, 21.88, 18.78, 17.86, 18.47, 15.13, 8.970000000000001, 4.55, 0, -0.02, -1.21, -5.6, -5, -2, -0.01, -0.01, 12.37, 19.29, 33.96, 39.96, 61.7, 53.03)
for ($i = 0; $i -lt $spotValues.Length; $i++) {
    $ws0.Cells.Item($i + 2, 5).Value = $spotValues[$i]
}

# --- New sheet: Gaz (inserted right after "Prix Spot") ---
$wsGaz = $wb.Worksheets.Add($null, $ws0)
$wsGaz.Name = "Gaz"

$gazHeaders = @("Date", "Contract", "Last", "High", "Low")
for ($i = 0; $i -lt $gazHeaders.Length; $i++) {
    $cell = $wsGaz.Cells.Item(1, $i + 1)
    $cell.Value = $gazHeaders[$i]
}
$ws0.Range("A1").Copy()
$wsGaz.Range("A1:E1").PasteSpecial(-4122)

$gazDateCell = $wsGaz.Cells.Item(2, 1)
$gazDateCell.NumberFormat = "@"
$gazDateCell.Value = "2025-06-17"
$gazDateCell.Style = "Normal"

$gazRow2Rest = @("PEG Day Ahead", "-", "-", "-")
for ($i = 0; $i -lt $gazRow2Rest.Length; $i++) {
    $wsGaz.Cells.Item(2, $i + 2).Value = $gazRow2Rest[$i]
}

# --- New sheet: CO2 (inserted right after "Gaz") ---
$wsCO2 = $wb.Worksheets.Add($null, $wsGaz)
$wsCO2.Name = "CO2"

$co2Headers = @("Date", "Last Price")
for ($i = 0; $i -lt $co2Headers.Length; $i++) {
    $co2Cell = $wsCO2.Cells.Item(1, $i + 1)
    $co2Cell.Value = $co2Headers[$i]
}
$ws0.Range("A1").Copy()
$wsCO2.Range("A1:B1").PasteSpecial(-4122)

$co2DateCell = $wsCO2.Cells.Item(2, 1)
$co2DateCell.NumberFormat = "@"
$co2DateCell.Value = "2025-06-17"
$co2DateCell.Style = "Normal"

$wsCO2.Cells.Item(2, 2).Value = "-"

# Restore original active sheet / selection (unchanged by the source diff)
$ws0.Activate()
$null = $ws0.Range("A1").Select()
